$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell references to their new text values.
# Values are written with a leading apostrophe to force Excel to
# store them as text (matching the inlineStr type used in the sheet)
# instead of letting it auto-convert numeric/percentage-looking
# strings into real numbers.
$updates = [ordered]@{
    "D2" = '291.46'
    "E2" = '-8.02%'
    "E3" = '-1.69%'
    "D4" = '5.012'
    "E4" = '-2.75%'
    "D5" = '0.07299'
    "E5" = '-4.46%'
    "D6" = '4.283'
    "E6" = '-0.76%'
    "D7" = '1.561'
    "E7" = '-7.31%'
    "D8" = '0.9221'
    "E8" = '-1.10%'
    "E9" = '-7.07%'
    "E10" = '-5.25%'
    "D11" = '0.08572'
    "E11" = '-5.11%'
    "D12" = '0.04183'
    "E12" = '0.45%'
    "E13" = '-0.31%'
    "D14" = '0.001262'
    "E14" = '-0.53%'
    "D15" = '0.005872'
    "E15" = '-0.78%'
    "D16" = '3.399'
    "E16" = '1.44%'
    "E18" = '-2.53%'
    "D19" = '7.840'
    "E19" = '-6.22%'
    "E20" = '2.54%'
    "D22" = '0.03853'
    "E22" = '-4.44%'
    "D23" = '0.001264'
    "E23" = '-0.66%'
    "D24" = '0.003788'
    "E24" = '-7.23%'
    "D25" = '0.0001281'
    "E25" = '0.42%'
    "D26" = '0.0003725'
    "E26" = '-95.03%'
    "D38" = '0.02306'
    "E38" = '-7.27%'
    "D39" = '0.04948'
    "E39" = '-5.48%'
    "E40" = '211.67%'
    "D41" = '0.007666'
    "E41" = '-1.55%'
    "D42" = '0.1275'
    "E42" = '-1.76%'
    "D43" = '0.007363'
    "E43" = '4.16%'
    "D44" = '0.007089'
    "E44" = '-13.96%'
    "D45" = '0.2908'
    "E45" = '-15.23%'
    "D46" = '0.00006414'
    "E46" = '-4.06%'
    "E47" = '-0.36%'
    "D48" = '0.01648'
    "E48" = '-92.47%'
    "E49" = '-0.41%'
    "D50" = '0.00002102'
    "E50" = '-0.36%'
    "D51" = '0.0002002'
    "E51" = '-0.36%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $updates[$cell]
    # Re-apply the default Normal style so the forced quote-prefix/text
    # formatting Excel applies for the apostrophe trick doesn't leave a
    # spurious style behind on the cell.
    $range.Style = "Normal"
}
